$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 6 and 7: swap Coin name + Link (FTXToken <-> KuCoinToken)
$ws.Range("B6").Value = "KuCoinToken"
$ws.Range("C6").Value = "https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs"
$ws.Range("B7").Value = "FTXToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"

# Columns D (Price), E (Volume 1h) and G (Hora) hold numeric-looking text
# (e.g. "330.20", "-0.61%", "13"). Force the Text number format on each
# touched cell right before writing its value, so Excel keeps the literal
# string instead of coercing it to a number (which would drop trailing
# zeros / % signs), without disturbing any other, untouched cells.

# Price (column D) updates
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "330.20"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "41.51"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.650"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.08332"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "8.788"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "2.020"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "4.532"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.992"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9218"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.1271"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.1962"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.09439"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.03870"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.1062"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.001307"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.006075"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.444"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "8.346"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.1375"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.2450"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04397"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001253"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004338"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0001197"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.02795"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.05518"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007760"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1426"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.008915"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.002239"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.01195"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006950"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.00000000748"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.003174"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00002095"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0001995"

# Volume(1h) (column E) updates
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "-0.61%"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "4.80%"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "-2.00%"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "3.26%"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "1.54%"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "0.32%"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "0.27%"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "1.32%"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "-0.17%"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "-0.97%"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "0.28%"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "2.52%"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "9.12%"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "0.94%"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "0.28%"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-5.10%"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "2.23%"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "1.43%"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "-4.54%"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "1.03%"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "-9.96%"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "-0.54%"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "-0.61%"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "-5.00%"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "-0.18%"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "10.88%"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "1.04%"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "3.63%"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "0.99%"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-10.11%"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "6.19%"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "5.38%"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "2.30%"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "-0.30%"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "4.74%"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "0.01%"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "-0.30%"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "-0.30%"

# Hora (column G) updates: 12 -> 13 for rows 2 through 51
for ($row = 2; $row -le 51; $row++) {
    $ws.Range("G$row").NumberFormat = "@"
    $ws.Range("G$row").Value = "13"
}
